# Daily attendance processing - 2025-11-04 10:25:53
# Reorders "Recorded By" email lists for several sessions, updates a few
# summary counters, and flips session B3/PHARMACOLOGY/2 (row 110) from
# "Pending" to "Not Recorded" (including its conditional-format coloring).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Recorded-By (column G) email re-orderings ----------------------------

$ws.Range("G2").Value   = "servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G3").Value   = "hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G7").Value   = "Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G13").Value  = "marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg"
$ws.Range("G18").Value  = "servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G19").Value  = "servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G23").Value  = "Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G29").Value  = "marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg"
$ws.Range("G34").Value  = "servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G35").Value  = "System, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G50").Value  = "servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G51").Value  = "System, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G66").Value  = "hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G67").Value  = "hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G82").Value  = "hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G83").Value  = "servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G98").Value  = "servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G103").Value = "Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G114").Value = "servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G119").Value = "Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# ---- Statistics counters ---------------------------------------------------

$ws.Range("L7").Value  = 8    # Missing Sessions (Class Statistics)
$ws.Range("L8").Value  = 95   # Pending Sessions (Class Statistics)
$ws.Range("P21").Value = 3    # Missing (Group Statistics, Year 2 / B3)
$ws.Range("Q21").Value = 11   # Pending (Group Statistics, Year 2 / B3)

# ---- Row 110 (Year 2 / B3 / PHARMACOLOGY / session 2): Pending -> Not Recorded

$ws.Range("I110").Value = "Not Recorded"

# Re-use the existing "Not Recorded" conditional formatting (fill/font) by
# copying the format from another already Not-Recorded row (109) onto row 110.
$src = $ws.Range("A109:I109")
$dst = $ws.Range("A110:I110")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false
